$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TK3 variables to netCDF")

# Row 6 of the "TK3 variables to netCDF" sheet gets a run of sequential index
# numbers typed across the (until now empty) tail of the header band, from
# column O through column BK. A couple of cells were left blank by the
# author (matches the commit message "Still missing some information"),
# and one value was mistyped (19 instead of 18), so the sequence below is
# reproduced exactly as entered rather than algorithmically regenerated.

$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 7
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 10
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 12
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 14
$ws.Range("AD6").Value = 15
$ws.Range("AE6").Value = 16
$ws.Range("AF6").Value = 17
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 20
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 22
$ws.Range("AK6").Value = 23
$ws.Range("AM6").Value = 24
$ws.Range("AN6").Value = 25
$ws.Range("AO6").Value = 26
$ws.Range("AP6").Value = 27
$ws.Range("AQ6").Value = 28
$ws.Range("AR6").Value = 29
$ws.Range("AS6").Value = 30
$ws.Range("AT6").Value = 31
$ws.Range("AU6").Value = 32
$ws.Range("AV6").Value = 33
$ws.Range("AW6").Value = 34
$ws.Range("AY6").Value = 35
$ws.Range("AZ6").Value = 36
$ws.Range("BA6").Value = 37
$ws.Range("BB6").Value = 38
$ws.Range("BC6").Value = 39
$ws.Range("BD6").Value = 40
$ws.Range("BE6").Value = 41
$ws.Range("BF6").Value = 42
$ws.Range("BG6").Value = 43
$ws.Range("BH6").Value = 44
$ws.Range("BI6").Value = 45
$ws.Range("BJ6").Value = 46
$ws.Range("BK6").Value = 47

# The user scrolled the frozen pane over to the right while entering this
# data and left the selection on the last cell they touched.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 32
$ws.Range("AZ12").Select()
